$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 84, shifting the existing rows (84-117) down to (85-118).
$ws.Rows("84").Insert()

# Fill in the new row 84 with the new weekly record. Columns that are constant
# across this whole "Tuna" dataset (A,B,C,E,F,G,H,I,J,K,Q,R,T) are copied from
# the row below (old row 84, now shifted to row 85), which Excel's Insert
# already carried the date-format style for (column D).
$ws.Range("A84").Value = 10
$ws.Range("B84").Value = "Vega Modelo de Temuco"
$ws.Range("C84").Value = "La Araucanía"
$ws.Range("D84").Value = 45229
$ws.Range("E84").Value = 9
$ws.Range("F84").Value = "Fruta"
$ws.Range("G84").Value = 100107
$ws.Range("H84").Value = "Otros"
$ws.Range("I84").Value = 100107011
$ws.Range("J84").Value = "Tuna"
$ws.Range("K84").Value = "Sin especificar"
$ws.Range("L84").Value = "Primera"
$ws.Range("M84").Value = 100
$ws.Range("N84").Value = 38000
$ws.Range("O84").Value = 38000
$ws.Range("P84").Value = 38000
$ws.Range("Q84").Value = '$/caja 16 kilos'
$ws.Range("R84").Value = "Provincia de Los Andes"
$ws.Range("S84").Value = 2375
$ws.Range("T84").Value = 16
